# This script applies the edit described by the diff:
# Two new rows of pricing data (Piña, Vega Monumental Concepción) are
# inserted right before the existing row 36, pushing all the existing
# data (previously rows 36-116) down by two rows (to rows 38-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 36, shifting rows 36:116 down to 38:118
$ws.Rows("36:37").Insert()

# --- New row 36 ---
$ws.Cells.Item(36, 1).Value  = 11
$ws.Cells.Item(36, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value  = "Bíobío"
$ws.Cells.Item(36, 4).Value  = 44519
$ws.Cells.Item(36, 5).Value  = 8
$ws.Cells.Item(36, 6).Value  = "Fruta"
$ws.Cells.Item(36, 7).Value  = 100108
$ws.Cells.Item(36, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(36, 9).Value  = 100108005
$ws.Cells.Item(36, 10).Value = "Piña"
$ws.Cells.Item(36, 11).Value = "Caramelo"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 100
$ws.Cells.Item(36, 14).Value = 18000
$ws.Cells.Item(36, 15).Value = 19000
$ws.Cells.Item(36, 16).Value = 18500
$ws.Cells.Item(36, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(36, 18).Value = "Ecuador"
$ws.Cells.Item(36, 19).Value = 1542
$ws.Cells.Item(36, 20).Value = 12

# --- New row 37 ---
$ws.Cells.Item(37, 1).Value  = 11
$ws.Cells.Item(37, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value  = "Bíobío"
$ws.Cells.Item(37, 4).Value  = 44519
$ws.Cells.Item(37, 5).Value  = 8
$ws.Cells.Item(37, 6).Value  = "Fruta"
$ws.Cells.Item(37, 7).Value  = 100108
$ws.Cells.Item(37, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value  = 100108005
$ws.Cells.Item(37, 10).Value = "Piña"
$ws.Cells.Item(37, 11).Value = "Caramelo"
$ws.Cells.Item(37, 12).Value = "Tercera"
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 18000
$ws.Cells.Item(37, 15).Value = 19000
$ws.Cells.Item(37, 16).Value = 18500
$ws.Cells.Item(37, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(37, 18).Value = "Ecuador"
$ws.Cells.Item(37, 19).Value = 1542
$ws.Cells.Item(37, 20).Value = 12
